$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 288. This shifts all existing rows
# 288..340 down to 290..342 (bringing their formatting/styles along),
# which already reproduces the target state for those rows.
$ws.Range("A288:A289").EntireRow.Insert()

# New row 288 (Primera) - newest week's data
$ws.Range("A288").Value = 8
$ws.Range("B288").Value = "Terminal La Palmera de La Serena"
$ws.Range("C288").Value = "Coquimbo"
$ws.Range("D288").Value = 44617
$ws.Range("E288").Value = 4
$ws.Range("F288").Value = 100112017
$ws.Range("G288").Value = "Apio"
$ws.Range("H288").Value = "Americana (o)"
$ws.Range("I288").Value = "Primera"
$ws.Range("J288").Value = 2200
$ws.Range("K288").Value = 8000
$ws.Range("L288").Value = 9000
$ws.Range("M288").Value = 8500
$ws.Range("N288").Value = "$/docena de matas"
$ws.Range("O288").Value = "Provincia del Elquí"
$ws.Range("P288").Value = 1417
$ws.Range("Q288").Value = 6
$ws.Range("R288").Value = "Hortaliza"

# New row 289 (Segunda) - newest week's data
$ws.Range("A289").Value = 8
$ws.Range("B289").Value = "Terminal La Palmera de La Serena"
$ws.Range("C289").Value = "Coquimbo"
$ws.Range("D289").Value = 44617
$ws.Range("E289").Value = 4
$ws.Range("F289").Value = 100112017
$ws.Range("G289").Value = "Apio"
$ws.Range("H289").Value = "Americana (o)"
$ws.Range("I289").Value = "Segunda"
$ws.Range("J289").Value = 1360
$ws.Range("K289").Value = 6000
$ws.Range("L289").Value = 7000
$ws.Range("M289").Value = 6500
$ws.Range("N289").Value = "$/docena de matas"
$ws.Range("O289").Value = "Provincia del Elquí"
$ws.Range("P289").Value = 1083
$ws.Range("Q289").Value = 6
$ws.Range("R289").Value = "Hortaliza"
